$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Wrong_Entity_NonEvent_as_Event"
$ws.Range("B2").Value = 63

$ws.Range("A3").Value = "Correct"
$ws.Range("B3").Value = 62

$ws.Range("A4").Value = "Wrong_Entity_Event_as_NonEvent"
$ws.Range("B4").Value = 40

$ws.Range("A5").Value = "Wrong_Tag_E_as_I"

$ws.Range("A6").Value = "Wrong_Tag_B_as_I"
